$wb = $excel.ActiveWorkbook

# Rename the sheets (matching the "Modification document du groupe" commit)
$wb.Worksheets.Item("UserStories").Name = "UserStories-Backlog"
$wb.Worksheets.Item("BurndownChart1").Name = "BurndownChart1-Sprint1"
$wb.Worksheets.Item("BurndownChart2").Name = "BurndownChart2-Sprint2"

# Renaming a sheet does not retarget existing chart series formulas, so
# update each chart's SERIES() formulas to point at the new sheet names.
$sheet1 = $wb.Worksheets.Item("BurndownChart1-Sprint1")
$chart1 = $sheet1.ChartObjects().Item(1).Chart
$chart1.SeriesCollection().Item(1).Formula = "=SERIES('BurndownChart1-Sprint1'!`$E`$12,'BurndownChart1-Sprint1'!`$I`$1:`$O`$1,'BurndownChart1-Sprint1'!`$I`$12:`$O`$12,1)"
$chart1.SeriesCollection().Item(2).Formula = "=SERIES('BurndownChart1-Sprint1'!`$E`$13,,'BurndownChart1-Sprint1'!`$I`$13:`$O`$13,2)"

$sheet2 = $wb.Worksheets.Item("BurndownChart2-Sprint2")
$chart2 = $sheet2.ChartObjects().Item(1).Chart
$chart2.SeriesCollection().Item(1).Formula = "=SERIES('BurndownChart2-Sprint2'!`$E`$12,'BurndownChart2-Sprint2'!`$I`$1:`$O`$1,'BurndownChart2-Sprint2'!`$I`$12:`$O`$12,1)"
$chart2.SeriesCollection().Item(2).Formula = "=SERIES('BurndownChart2-Sprint2'!`$E`$13,,'BurndownChart2-Sprint2'!`$I`$13:`$O`$13,2)"
